$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9850.272000000001
$ws.Range("I9").Value = 14385.857
$ws.Range("K9").Value = 14385.857
$ws.Range("M9").Value = -14216.857
$ws.Range("H19").Value = 1175.6154
$ws.Range("I19").Value = 1128.3
$ws.Range("K19").Value = 1128.3
$ws.Range("M19").Value = -953.3
$ws.Range("H42").Value = 203.85715
$ws.Range("I42").Value = 42.333332
$ws.Range("J42").Value = 325
$ws.Range("K42").Value = 126.999996
$ws.Range("L42").Value = 975
$ws.Range("M42").Value = 103.000004
$ws.Range("N42").Value = -1435
$ws.Range("H53").Value = 591.2632
$ws.Range("I53").Value = 418.3
$ws.Range("J53").Value = 783.44446
$ws.Range("K53").Value = 418.3
$ws.Range("L53").Value = 783.44446
$ws.Range("M53").Value = 218.7
$ws.Range("N53").Value = -2057.44446
$ws.Range("H80").Value = 2143.182
$ws.Range("J80").Value = 2819.7856
$ws.Range("L80").Value = 8459.356800000001
$ws.Range("N80").Value = -10455.3568
$ws.Range("H83").Value = 2143.182
$ws.Range("J83").Value = 2819.7856
$ws.Range("L83").Value = 25378.0704
$ws.Range("N83").Value = -35362.0704
$ws.Range("H116").Value = 7353.273
$ws.Range("I116").Value = 6977.5
$ws.Range("J116").Value = 11111
$ws.Range("K116").Value = 6977.5
$ws.Range("L116").Value = 11111
$ws.Range("M116").Value = -3535.5
$ws.Range("N116").Value = -17995
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("N130").Value = 0
$ws.Range("H131").Value = 4222.25
$ws.Range("J131").Value = 6900
$ws.Range("L131").Value = 20700
$ws.Range("N131").Value = -30780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 249.83333
$ws.Range("J4").Value = 224.75
$ws.Range("L4").Value = 224.75
$ws.Range("M4").Value = -184
$ws.Range("N4").Value = -456.75
$ws.Range("H63").Value = 1927.2727
$ws.Range("I63").Value = 1927.2727
$ws.Range("K63").Value = 1927.2727
$ws.Range("M63").Value = -1241.2727
$ws.Range("H66").Value = 1927.2727
$ws.Range("I66").Value = 1927.2727
$ws.Range("K66").Value = 9636.363499999999
$ws.Range("M66").Value = -6204.363499999999
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("N128").Value = 0
$ws.Range("H130").Value = 63971.6
$ws.Range("J130").Value = 63971.6
$ws.Range("L130").Value = 63971.6
$ws.Range("N130").Value = -74011.60000000001
$ws.Range("H131").Value = 49999
$ws.Range("J131").Value = 49999
$ws.Range("L131").Value = 49999
$ws.Range("N131").Value = -60079
$ws.Range("H132").Value = 35567.453
$ws.Range("I132").Value = 3853.2432
$ws.Range("K132").Value = 11559.7296
$ws.Range("M132").Value = -9029.729599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 621.55554
$ws.Range("I22").Value = 427.42856
$ws.Range("K22").Value = 427.42856
$ws.Range("M22").Value = -254.42856
$ws.Range("H107").Value = 3295.9
$ws.Range("I107").Value = 3295.9
$ws.Range("K107").Value = 3295.9
$ws.Range("M107").Value = -1375.9
$ws.Range("H123").Value = 66000
$ws.Range("J123").Value = 66000
$ws.Range("L123").Value = 66000
$ws.Range("N123").Value = -75800
$ws.Range("H124").Value = 101885
$ws.Range("J124").Value = 101885
$ws.Range("L124").Value = 101885
$ws.Range("N124").Value = -111705
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").Value = 0
$ws.Range("H126").Value = 50000
$ws.Range("J126").Value = 50000
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("N127").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("N129").Value = 0
$ws.Range("H130").Value = 80562
$ws.Range("J130").Value = 80562
$ws.Range("L130").Value = 80562
$ws.Range("N130").Value = -90602
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("N131").Value = 0
$ws.Range("H134").Value = 3653.6667
$ws.Range("I134").Value = 12000
$ws.Range("J134").Value = 1984.4
$ws.Range("K134").Value = 36000
$ws.Range("L134").Value = 5953.200000000001
$ws.Range("M134").Value = -33465
$ws.Range("N134").Value = -11023.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 224.11765
$ws.Range("I7").Value = 58
$ws.Range("J7").Value = 340.4
$ws.Range("K7").Value = 58
$ws.Range("L7").Value = 340.4
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = -566.4
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("H22").Value = 5122.6816
$ws.Range("I22").Value = 5833.278
$ws.Range("K22").Value = 5833.278
$ws.Range("M22").Value = -5483.278
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("H62").Value = 3487.875
$ws.Range("J62").Value = 3500.6667
$ws.Range("L62").Value = 3500.6667
$ws.Range("N62").Value = -4748.6667
$ws.Range("H65").Value = 3487.875
$ws.Range("J65").Value = 3500.6667
$ws.Range("L65").Value = 17503.3335
$ws.Range("N65").Value = -23743.3335
$ws.Range("H122").Value = 5849559.5
$ws.Range("I122").Value = 1494
$ws.Range("K122").Value = 4482
$ws.Range("M122").Value = -2032
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("N128").Value = 0
$ws.Range("H131").Value = 49599.75
$ws.Range("J131").Value = 49599.75
$ws.Range("L131").Value = 49599.75
$ws.Range("N131").Value = -59679.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 586.8570999999999
$ws.Range("I23").Value = 466.8
$ws.Range("J23").Value = 653.55554
$ws.Range("K23").Value = 1400.4
$ws.Range("L23").Value = 1960.66662
$ws.Range("M23").Value = -1165.4
$ws.Range("N23").Value = -2430.66662
$ws.Range("H68").Value = 1000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("H71").Value = 1000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("H80").Value = 2951.5
$ws.Range("J80").Value = 2951.5
$ws.Range("L80").Value = 8854.5
$ws.Range("N80").Value = -10726.5
$ws.Range("H83").Value = 2951.5
$ws.Range("J83").Value = 2951.5
$ws.Range("L83").Value = 26563.5
$ws.Range("N83").Value = -35923.5
$ws.Range("H106").Value = 3352.6667
$ws.Range("J106").Value = 4029
$ws.Range("L106").Value = 12087
$ws.Range("N106").Value = -13979
$ws.Range("H121").Value = 1549.5555
$ws.Range("J121").Value = 1278.1428
$ws.Range("L121").Value = 3834.4284
$ws.Range("N121").Value = -6454.428400000001
$ws.Range("H131").Value = 27779220
$ws.Range("J131").Value = 1686.5
$ws.Range("L131").Value = 5059.5
$ws.Range("N131").Value = -15139.5
$ws.Range("H133").Value = 2427.5715
$ws.Range("I133").Value = 1482.5
$ws.Range("K133").Value = 4447.5
$ws.Range("M133").Value = 612.5
$ws.Range("H134").Value = 1038.8334
$ws.Range("I134").Value = 921.5294
$ws.Range("J134").Value = 3033
$ws.Range("K134").Value = 2764.5882
$ws.Range("L134").Value = 9099
$ws.Range("M134").Value = 2305.4118
$ws.Range("N134").Value = -19239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 10381.4
$ws.Range("I102").Value = 1726.75
$ws.Range("K102").Value = 1726.75
$ws.Range("M102").Value = -104.75
$ws.Range("H136").Value = 32317.857
$ws.Range("J136").Value = 32317.857
$ws.Range("L136").Value = 96953.571
$ws.Range("N136").Value = -102053.571
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3226.1177
$ws.Range("I61").Value = 2722.9333
$ws.Range("K61").Value = 2722.9333
$ws.Range("M61").Value = -2520.9333
$ws.Range("H113").Value = 3226.1177
$ws.Range("I113").Value = 2722.9333
$ws.Range("K113").Value = 2722.9333
$ws.Range("M113").Value = -552.9333000000001
$ws.Range("H133").Value = 58330
$ws.Range("J133").Value = 58330
$ws.Range("L133").Value = 58330
$ws.Range("N133").Value = -63390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("N129").Value = 0
$ws.Range("H130").Value = 47024.5
$ws.Range("J130").Value = 47024.5
$ws.Range("L130").Value = 47024.5
$ws.Range("N130").Value = -57064.5
